$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.097.35'
$ws.Range('E2').Value = '  -3.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.514.63'
$ws.Range('E3').Value = '  -4.74%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.79'
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.04'
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.507.90'
$ws.Range('E8').Value = '  -4.68%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -5.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.77'
$ws.Range('E11').Value = '  +7.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.600'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.34'
$ws.Range('E13').Value = '  -5.43%  '
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '677.73'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.082.24'
$ws.Range('E16').Value = '  -4.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.78'
$ws.Range('E17').Value = '  -2.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.515.63'
$ws.Range('E18').Value = '  -4.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.101.10'
$ws.Range('E19').Value = '  -3.73%  '
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.59'
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.26'
$ws.Range('E22').Value = '  -3.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.908'
$ws.Range('E23').Value = '  -3.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.18'
$ws.Range('E24').Value = '  -9.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.34'
$ws.Range('E25').Value = '  -5.53%  '
$ws.Range('E26').Value = '  -4.22%  '
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -5.97%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -6.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.02'
$ws.Range('E31').Value = '  -7.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.80'
$ws.Range('E32').Value = '  -4.90%  '
$ws.Range('E33').Value = '  -7.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.39'
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('E35').Value = '  -5.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '580.53'
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.60'
$ws.Range('E37').Value = '  -14.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.94'
$ws.Range('E38').Value = '  -3.31%  '
$ws.Range('E39').Value = '  -3.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '57.42'
$ws.Range('E40').Value = '  -3.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -3.78%  '
$ws.Range('E43').Value = '  -5.00%  '
$ws.Range('E44').Value = '  -6.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.435.39'
$ws.Range('E45').Value = '  -9.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '33.56'
$ws.Range('E46').Value = '  -5.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₃0711'
$ws.Range('E47').Value = '  -8.98%  '
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '132.14'
$ws.Range('E51').Value = '  -1.98%  '
